$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 16
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = "2025-03-23 21:55:01"
$ws.Range("E2").Value = 6000.06
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "O254"

# Update row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 9
$ws.Range("D3").Value = "2025-03-23 21:59:22"
$ws.Range("E3").Value = 8000.08
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "O300"

# Delete rows 4-8 (entire rows)
$ws.Range("A4:G8").Delete()
